$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price (D) column updates ---

# Some D-column values are valid numeric literals and Excel would normally
# auto-convert them to numbers on assignment. The source data keeps every
# Price cell as plain text, so force Text format before assigning, then
# restore the default "Normal" style (keeps the value as text, no format
# residue on the cell).
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '155.04'
$ws.Range("D6").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '6.68'
$ws.Range("D9").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.130'
$ws.Range("D12").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '355.26'
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '10.58'
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '1.00'
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '61.02'
$ws.Range("D23").Style = "Normal"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.41'
$ws.Range("D29").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '6.30'
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.42'
$ws.Range("D32").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '147.98'
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.21'
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.934'
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.21'
$ws.Range("D37").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.81'
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '36.48'
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '291.48'
$ws.Range("D42").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.623'
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0561'
$ws.Range("D45").Style = "Normal"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '19.62'
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0238'
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.26'
$ws.Range("D51").Style = "Normal"

# D-column values that are not valid numeric literals (contain multiple
# separators or special glyphs) -- Excel keeps these as text automatically.
$ws.Range("D2").Value = '61.012.56'
$ws.Range("D3").Value = '2.605.24'
$ws.Range("D13").Value = '3.061.12'
$ws.Range("D14").Value = '61.026.61'
$ws.Range("D17").Value = '2.607.53'
$ws.Range("D26").Value = '2.716.02'
$ws.Range("D28").Value = '0.0₃0849'

# --- Volume(1h) (E) column updates ---
# These strings include "%" / leading "+"/"-" and padding spaces, so they
# are never misread as numbers and stay as plain text automatically.
$ws.Range("E2").Value = '  +1.43%  '
$ws.Range("E3").Value = '  +1.36%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("E5").Value = '  +3.59%  '
$ws.Range("E6").Value = '  +1.61%  '
$ws.Range("E7").Value = '  -0.08%  '
$ws.Range("E8").Value = '  +2.59%  '
$ws.Range("E9").Value = '  +2.06%  '
$ws.Range("E10").Value = '  +1.95%  '
$ws.Range("E11").Value = '  +0.05%  '
$ws.Range("E12").Value = '  +1.28%  '
$ws.Range("E13").Value = '  +1.39%  '
$ws.Range("E14").Value = '  +1.45%  '
$ws.Range("E15").Value = '  +1.14%  '
$ws.Range("E16").Value = '  +1.53%  '
$ws.Range("E17").Value = '  +1.40%  '
$ws.Range("E18").Value = '  +0.05%  '
$ws.Range("E19").Value = '  +2.97%  '
$ws.Range("E20").Value = '  +1.82%  '
$ws.Range("E21").Value = '  +2.28%  '
$ws.Range("E22").Value = '  +0.15%  '
$ws.Range("E23").Value = '  +2.25%  '
$ws.Range("E24").Value = '  +2.14%  '
$ws.Range("E25").Value = '  +0.76%  '
$ws.Range("E26").Value = '  +1.11%  '
$ws.Range("E28").Value = '  +1.05%  '
$ws.Range("E29").Value = '  +0.50%  '
$ws.Range("E30").Value = '  -0.07%  '
$ws.Range("E31").Value = '  +10.64%  '
$ws.Range("E32").Value = '  +0.76%  '
$ws.Range("E33").Value = '  +3.07%  '
$ws.Range("E34").Value = '  -3.40%  '
$ws.Range("E35").Value = '  +5.83%  '
$ws.Range("E36").Value = '  +9.24%  '
$ws.Range("E37").Value = '  +2.14%  '
$ws.Range("E38").Value = '  +2.82%  '
$ws.Range("E39").Value = '  +2.68%  '
$ws.Range("E40").Value = '  +2.20%  '
$ws.Range("E41").Value = '  +1.22%  '
$ws.Range("E42").Value = '  -0.28%  '
$ws.Range("E43").Value = '  +2.47%  '
$ws.Range("E44").Value = '  +1.10%  '
$ws.Range("E45").Value = '  +1.37%  '
$ws.Range("E46").Value = '  -0.04%  '
$ws.Range("E47").Value = '  +4.50%  '
$ws.Range("E48").Value = '  -0.19%  '
$ws.Range("E49").Value = '  +2.43%  '
$ws.Range("E50").Value = '  +0.14%  '
$ws.Range("E51").Value = '  +9.22%  '
